$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "trainingimages/13_kopopi"
$ws.Range("B2").Value = "pngimages/13_toast.png"
$ws.Range("C2").Value = "trainingimages/02_pitito"
$ws.Range("D2").Value = "pngimages/02_pallet.png"
$ws.Range("E2").Value = 0.5
$ws.Range("F2").Value = -0.5

# Row 3
$ws.Range("A3").Value = "trainingimages/27_pakapa"
$ws.Range("B3").Value = "pngimages/27_kiwi.png"
$ws.Range("C3").Value = "trainingimages/09_tipata"
$ws.Range("D3").Value = "pngimages/09_plane.png"

# Row 4
$ws.Range("A4").Value = "trainingimages/07_pitapi"
$ws.Range("B4").Value = "pngimages/07_suitcase.png"
$ws.Range("C4").Value = "trainingimages/24_takopa"
$ws.Range("D4").Value = "pngimages/24_banana.png"
